# Updates the cryptos price/volume table with refreshed market data.
# D-column price cells are numeric-looking strings (e.g. "38.750.64",
# "0.999") that must stay plain text, matching the workbook's existing
# inlineStr cells. A leading apostrophe forces Excel to keep them as text
# instead of auto-converting to a Number, and resetting the cell Style
# back to "Normal" afterwards clears the transient quote-prefix formatting
# so only the cell value changes (no stray style/format diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.750.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "'2.103.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'228.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'62.16"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.389"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").Value = "'0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "'15.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.63%  "
$ws.Range("D13").Value = "'2.413.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "'22.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").Value = "'0.807"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("D16").Value = "'5.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "'2.107.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "'38.777.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "'71.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").Value = "'6.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").Value = "'0.0₃0841"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").Value = "'227.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").Value = "'2.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'172.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").Value = "'9.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.65%  "
$ws.Range("E28").Value = "  +5.28%  "
$ws.Range("E29").Value = "  +4.54%  "
$ws.Range("D30").Value = "'19.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("D31").Value = "'2.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.19%  "
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").Value = "'4.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("E35").Value = "  +8.05%  "
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").Value = "'2.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("D38").Value = "'3.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "'18.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").Value = "'103.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("E42").Value = "  +4.03%  "
$ws.Range("D43").Value = "'1.537.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'1.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.57%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'7.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.30%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "'2.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "'2.300.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.17%  "
